$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update resolution (ppi) values for "art" and "type" alias rows to 300
$ws.Range("E8").Value = 300
$ws.Range("E9").Value = 300

# Update the active selection to match the author's last selection
$ws.Range("F9").Select()
